# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) figures across the four sheets, and
# inserts a newly-scraped duplicate listing for the 9/1 "音阅派国漫演唱会"
# event into the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) - F column refresh
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1173
$ws1.Range("F4").Value = 12584
$ws1.Range("F5").Value = 703
$ws1.Range("F10").Value = 318
$ws1.Range("F11").Value = 1845
$ws1.Range("F12").Value = 38
$ws1.Range("F17").Value = 328
$ws1.Range("F19").Value = 284
$ws1.Range("F20").Value = 117
$ws1.Range("F21").Value = 118
$ws1.Range("F23").Value = 206
$ws1.Range("F24").Value = 235
$ws1.Range("F25").Value = 1252
$ws1.Range("F26").Value = 59

# ---------------------------------------------------------------------
# Sheet "演出" (performances) - F column refresh
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 273
$ws2.Range("F6").Value = 127
$ws2.Range("F9").Value = 68
$ws2.Range("F10").Value = 349
$ws2.Range("F15").Value = 10

# ---------------------------------------------------------------------
# Sheet "本地生活" (local life) - F column refresh
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 860

# ---------------------------------------------------------------------
# Sheet "全部类型" (combined) - F column refresh + new row insertion
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

# Insert a new row 27 that duplicates the existing 9/1 "音阅派" listing
# (currently sitting in row 26), pushing everything below down by one.
$ws4.Rows.Item(27).Insert()

# Copy row 26's formatting (column A carries the bordered/centered style)
# onto the new row 27's A cell so the inserted row looks native.
$ws4.Range("A26").Copy()
$ws4.Range("A27").PasteSpecial(-4122)

$ws4.Range("A27").Value = 26
$ws4.Range("B27").Value = "'2024-09-01"
$ws4.Range("C27").Value = "广州·音阅派国漫演唱会-《狐妖小红娘》《一人之下》领衔国漫原声音乐现场"
$ws4.Range("D27").Value = "东风中路259号 广州中山纪念堂"
$ws4.Range("E27").Value = "2024.09.01 19:30-09.01 21:00"
$ws4.Range("F27").Value = 68
$ws4.Range("G27").Value = 180
$ws4.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=89794"
$ws4.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202407/ehol1aeZ1721810539643.jpeg"

# Re-sequence column A (the plain ordinal index) for every row pushed down
# by the insert - these are literal values, not formulas, so Excel does not
# renumber them automatically.
For ($r = 28; $r -le 48; $r++) {
    $ws4.Range("A$r").Value = $r - 1
}

# Refresh the F column ("想去人数") to match the latest source-sheet figures,
# at the new (post-insert) row positions.
$ws4.Range("F2").Value = 860
$ws4.Range("F6").Value = 1173
$ws4.Range("F7").Value = 12584
$ws4.Range("F8").Value = 273
$ws4.Range("F9").Value = 703
$ws4.Range("F14").Value = 318
$ws4.Range("F15").Value = 1845
$ws4.Range("F16").Value = 38
$ws4.Range("F21").Value = 127
$ws4.Range("F22").Value = 127
$ws4.Range("F26").Value = 68
$ws4.Range("F28").Value = 349
$ws4.Range("F29").Value = 328
$ws4.Range("F32").Value = 284
$ws4.Range("F33").Value = 117
$ws4.Range("F34").Value = 118
$ws4.Range("F37").Value = 206
$ws4.Range("F40").Value = 235
$ws4.Range("F41").Value = 1252
$ws4.Range("F42").Value = 10
$ws4.Range("F43").Value = 59
